# Update odds values in Sheet1 to reflect the latest Betfair Back/Lay data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 1.65
$ws.Range("G2").Value = 1.71
$ws.Range("H2").Value = 5.2
$ws.Range("I2").Value = 5.6
$ws.Range("J2").Value = 4.4
$ws.Range("R2").Value = 1.52
$ws.Range("W2").Value = 2.38
$ws.Range("AK2").Value = 16.5

# Row 4
$ws.Range("F4").Value = 1.99
$ws.Range("H4").Value = 2.84
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 7.2

# Row 6
$ws.Range("H6").Value = 5
$ws.Range("K6").Value = 9.6

# Row 7
$ws.Range("F7").Value = 2.28
$ws.Range("J7").Value = 3.05
$ws.Range("K7").Value = 5.6
$ws.Range("P7").Value = 1.71
$ws.Range("Q7").Value = 1.91

# Row 8
$ws.Range("F8").Value = 1.52
$ws.Range("G8").Value = 1.85
$ws.Range("Q8").Value = 2.5

# Row 10
$ws.Range("H10").Value = 32
$ws.Range("K10").Value = 13
